# Atualização de bases das ligas, do dia: 09-04-2024 às 22:40
#
# The source data contained several duplicate-looking fixture rows whose
# home/away teams (and the rest of the odds columns) had been mixed up
# between two (or, in one case, three) rows. This script swaps the B:AC
# content of each affected pair/group of rows back into the correct rows,
# and also applies a handful of standalone odds corrections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

function Get-RowValues($sheet, [int]$row) {
    $vals = New-Object 'object[]' ($lastCol - $firstCol + 1)
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals[$c - $firstCol] = $sheet.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($sheet, [int]$row, $vals) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $sheet.Cells.Item($row, $c).Value = $vals[$c - $firstCol]
    }
}

function Swap-Rows($sheet, [int]$rowA, [int]$rowB) {
    $a = Get-RowValues $sheet $rowA
    $b = Get-RowValues $sheet $rowB
    Set-RowValues $sheet $rowA $b
    Set-RowValues $sheet $rowB $a
}

# Simple pairwise swaps (home/away + odds got attributed to the wrong one
# of two rows sharing the same date).
$pairs = @(
    @(34, 35),
    @(164, 165),
    @(175, 176),
    @(181, 182),
    @(183, 184),
    @(187, 188),
    @(190, 191)
)

foreach ($pair in $pairs) {
    Swap-Rows $ws $pair[0] $pair[1]
}

# Rows 241 / 243 / 244 form a 3-way rotation rather than a simple swap:
#   new(241) = old(244)
#   new(243) = old(241)
#   new(244) = old(243)
$v241 = Get-RowValues $ws 241
$v243 = Get-RowValues $ws 243
$v244 = Get-RowValues $ws 244

Set-RowValues $ws 241 $v244
Set-RowValues $ws 243 $v241
Set-RowValues $ws 244 $v243

# A handful of standalone odds corrections (no row/id changes involved).
$ws.Cells.Item(259, 18).Value = 1.825   # R259
$ws.Cells.Item(259, 19).Value = 2.025   # S259

$ws.Cells.Item(260, 21).Value = 1.9     # U260
$ws.Cells.Item(260, 22).Value = 1.95    # V260

$ws.Cells.Item(263, 18).Value = 2.025   # R263
$ws.Cells.Item(263, 19).Value = 1.825   # S263

$ws.Cells.Item(265, 18).Value = 2       # R265
$ws.Cells.Item(265, 19).Value = 1.85    # S265
